$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: merge the two runs "Qual é o paralelismo médio" + "?" (the
# first occurrence in the document, right after "Grafo da direita: 4.")
# into a single run "Qual é o paralelismo médio?". Using wdReplaceOne (1)
# so only the first (already-split) occurrence is touched; the second
# occurrence later in the doc is already a single run and must stay
# untouched.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Qual é o paralelismo médio?", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Qual é o paralelismo médio?", 1) | Out-Null

# ---------------------------------------------------------------------
# Edit 2: add a collapsed "_GoBack" bookmark right after the run
# "64/34 =~ 1.882" (end of that paragraph's text, before the paragraph
# mark). A directly-collapsed Range placed exactly at that boundary
# does not anchor reliably, so we insert a temporary character there,
# anchor the bookmark around it, then delete the temporary character -
# the bookmark collapses back to the correct position.
# ---------------------------------------------------------------------
$p64 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*64/34 =~ 1.882*") {
        $p64 = $p
    }
}
$r64 = $p64.Range.Duplicate
$r64.MoveEnd(1, -1)
$r64.Collapse(0)
$r64.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $r64) | Out-Null
$r64.Text = ""

# ---------------------------------------------------------------------
# Edit 3: make the "Slide nº 15" paragraph (the one that follows the
# empty paragraph right after "64/34 =~ 1.882") bold - both the run and
# the paragraph mark formatting.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Slide nº 15*") {
        $p.Range.Bold = 1
    }
}

# ---------------------------------------------------------------------
# Edit 4: remove the "_GoBack" bookmark that originally sat right after
# the run "4+3n." (it has effectively moved earlier in the document, to
# right after "64/34 =~ 1.882", per edit 2 above).
# ---------------------------------------------------------------------
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# ---------------------------------------------------------------------
# Edit 5: merge the two runs "Qual é o grau máximo de concorrência" +
# "?" into a single run "Qual é o grau máximo de concorrência?". This
# must only touch the occurrence right after "4+3n." - the other
# occurrence earlier in the document is already a single run. Scope the
# Find to the range starting right after "4+3n." to target the correct
# one.
# ---------------------------------------------------------------------
$fullText = $d.Content.Text
$afterIdx = $fullText.IndexOf("4+3n.")
$scopeStart = $afterIdx + 5
$scope = $d.Range($scopeStart, $d.Content.End)
$scope.Find.Execute("Qual é o grau máximo de concorrência?", $false, $false, $false, $false, $false,
                     $true, 1, $false, "Qual é o grau máximo de concorrência?", 1) | Out-Null

Write-Output "edits applied"
